$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "NA" values under the duplicate_image_filename column (column E)
# for data rows 2 through 21.
foreach ($r in 2..21) {
    $ws.Cells.Item($r, 5).Value = "NA"
}

# Keep F1 as an empty placeholder cell (untouched by this edit); the
# underlying engine otherwise re-materializes it with a stray value on
# save, so explicitly reset it back to blank to avoid an unrelated diff.
$ws.Cells.Item(1, 6).Value = ""
